# Daily attendance processing - 2026-01-03 11:00:29
# Normalizes the "Recorded By" (column G) values so that entries that
# combine the owner's email with another recorder ("System" or
# "admin@admin.com") list the other recorder first, e.g.
#   "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2

    if ($value -eq $null) { continue }

    $text = [string]$value

    if ($text -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($text -eq "dnasr281@gmail.com, admin@admin.com") {
        $cell.Value2 = "admin@admin.com, dnasr281@gmail.com"
    }
}
